# Applies italic/bold run-property normalization across all slides:
#  - Title runs (44pt, bold): add explicit i="0" (not italic)
#  - Body bullet runs (28pt): add explicit b="0" i="0" (not bold, not italic)
$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)

    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $shape = $s.Shapes.Item($j)
        if (-not $shape.HasTextFrame) { continue }

        $tf = $shape.TextFrame
        $tr = $tf.TextRange
        $paraCount = $tr.Paragraphs().Count

        for ($k = 1; $k -le $paraCount; $k++) {
            $para = $tr.Paragraphs($k)

            $isTitle = $false
            if ($shape.Type -eq 14) {
                # msoPlaceholder: 1 = ppPlaceholderTitle / ppPlaceholderCenterTitle-ish title
                if ($shape.PlaceholderFormat.Type -eq 1 -or $shape.PlaceholderFormat.Type -eq 13) {
                    $isTitle = $true
                }
            }

            if ($isTitle) {
                # Title placeholder: keep bold, explicitly mark not-italic.
                $para.Font.Italic = 0
            } else {
                # Body placeholder: explicitly mark not-bold, not-italic.
                $para.Font.Bold = 0
                $para.Font.Italic = 0
            }
        }
    }
}
